$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.748.34"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "2.338.29"
$ws.Range("E3").Value = "  -0.85%  "

$ws.Range("E4").Value = "  +0.02%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.55"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -1.64%  "

$ws.Range("E6").Value = "  -4.13%  "

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.79"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  -7.41%  "

$ws.Range("E8").Value = "  +0.00%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.597"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -4.97%  "

$ws.Range("E10").Value = "  -3.37%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.80"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +0.75%  "

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.38"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  -4.40%  "

$ws.Range("E13").Value = "  -0.66%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.13"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -6.44%  "

$ws.Range("D15").Value = "2.684.86"
$ws.Range("E15").Value = "  -0.71%  "

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.02"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -6.05%  "

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.898"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -3.61%  "

$ws.Range("D18").Value = "2.345.29"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D19").Value = "43.681.22"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("E20").Value = "  -2.80%  "

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "77.57"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("E22").Value = "  -1.65%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.99"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -3.32%  "

$ws.Range("E24").Value = "  +0.07%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.89"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +5.50%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.71"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +2.47%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.48"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -2.69%  "

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.29"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -6.65%  "

$ws.Range("E29").Value = "  -1.31%  "

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "176.11"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +0.62%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.09"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -4.11%  "

$ws.Range("E32").Value = "  -2.35%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.132"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -2.76%  "

$ws.Range("E34").Value = "  -3.32%  "

$ws.Range("E35").Value = "  -5.63%  "

$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.31"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -1.87%  "

$ws.Range("E37").Value = "  -2.10%  "

$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.34"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -2.08%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.35"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -3.47%  "

$ws.Range("B40").Value = "FTXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.75"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +27.72%  "

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0269"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -3.42%  "

$ws.Range("E42").Value = "  +17.77%  "

$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.107"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +3.25%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.10"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +0.74%  "

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.59"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -4.83%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.195"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -3.67%  "

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.22"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -3.73%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.41"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -5.33%  "

$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.92"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +2.89%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "97.57"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -4.22%  "
